$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(305487936, Avihai  Kipnis: 6,0)"
$ws.Range("B1").Value = "(313227928, Aviv  Levi: 1,8)"
$ws.Range("C1").Value = "(205807308, Sariel  Basis: 5,4)"
$ws.Range("D1").Value = "(315891549, Raz  Halaby: 0,4)"
$ws.Range("E1").Value = "(315060103, Dan  Mshelh: 2,6)"
$ws.Range("F1").Value = "(313925141, Elad   Amer: 2,-3)"
$ws.Range("G1").Value = "(326598423, Ron Cohen: 3,8)"

$ws.Range("A3").Value = "cost: 415.9898314738982"
$ws.Range("A4").Value = "time: 55.855690210556894"
